# Remove the "Tenure Years" column from the Monthly Performance sheet.
# It sits between "Reporting Head" (N1) and "HR Comments" (O1/P1),
# so shift "HR Comments" left into O1 and clear the now-unused P1 cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = $ws.Range("P1").Value()
$ws.Range("P1").ClearContents()

$ws.Range("O1").Select()
